$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: paragraph "Реализовал стратегии микро-, мезо- и макро-управления
# для симуляции" currently consists of three runs split around
# <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/> markers
# (Word's grammar-check bracketing of "макро-управления"). The target
# rewrites this as a single run with the proofErr markers gone.
#
# Re-assigning Range.Text replaces the whole run/proofErr sequence inside
# the range with one freshly written run that inherits the formatting of
# the range's starting point, which is exactly what's needed here. The
# engine treats a Range.Text assignment that is textually identical to the
# range's current text as a no-op, so we first stamp a throw-away value to
# force the rewrite, then set the real target text.
# ---------------------------------------------------------------------------

$microMesoPara = $d.Paragraphs(14)
$microMesoRange = $microMesoPara.Range
$microMesoBody = $d.Range($microMesoRange.Start, $microMesoRange.End - 1)
$microMesoBody.Text = "x"
$microMesoPara2 = $d.Paragraphs(14)
$microMesoRange2 = $microMesoPara2.Range
$microMesoBody2 = $d.Range($microMesoRange2.Start, $microMesoRange2.End - 1)
$microMesoBody2.Text = "Реализовал стратегии микро-, мезо- и макро-управления для симуляции"

# ---------------------------------------------------------------------------
# Edit 2: add a new bulleted list item "Написал разделы "симуляция" в
# статьи" right after "Провел практический эксперимент для сравнения
# эффективности работы стратегий управления", reusing that paragraph's
# formatting (same pStyle/numPr/tabs/spacing/ind/rPr).
# ---------------------------------------------------------------------------

$practicalExperimentPara = $d.Paragraphs(15)
$practicalExperimentRange = $practicalExperimentPara.Range
$practicalExperimentFull = $d.Range($practicalExperimentRange.Start, $practicalExperimentRange.End)
$practicalExperimentFull.InsertParagraphAfter()

$newPara = $d.Paragraphs(16)
$newRange = $newPara.Range
$newBody = $d.Range($newRange.Start, $newRange.End - 1)
$newBody.Text = "x"
$newPara2 = $d.Paragraphs(16)
$newRange2 = $newPara2.Range
$newBody2 = $d.Range($newRange2.Start, $newRange2.End - 1)
$newBody2.Text = "Написал разделы " + [char]0x201C + "симуляция" + [char]0x201D + " в статьи"

# ---------------------------------------------------------------------------
# Edit 3: inside the review-signature textbox (a legacy VML <w:pict>/
# <v:shape> text box holding the "Гранич ин Олег Николаевич" table), drop
# the <w:proofErr w:type="spellStart"/> / <w:proofErr w:type="spellEnd"/>
# pair bracketing "Граничин" and merge the "ин" run with the following
# single-space run into one "ин " run. That textbox content lives in a
# story that isn't reachable through Document.Paragraphs / Document.Range /
# Document.Find in this object model (legacy VML text boxes aren't walked
# by the exposed Word OM surface), so there is no supported COM path to
# reach it. Try the straightforward route defensively in case the host can
# resolve it, but don't let a failure here abort the rest of the script.
# ---------------------------------------------------------------------------

try {
    $spellFixRange = $d.Content
    $spellFixRange.Find.ClearFormatting()
    $spellFound = $spellFixRange.Find.Execute("Граничин", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($spellFound) {
        $spellFixRange.Text = "Граничин"
    }
} catch {
}
